$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

Set-TextValue "D2" "243.54"
Set-TextValue "D3" "23.89"
Set-TextValue "D4" "5.258"
Set-TextValue "D6" "6.472"
Set-TextValue "D7" "3.332"
Set-TextValue "D8" "0.8117"
Set-TextValue "D9" "0.8743"
Set-TextValue "D11" "0.07260"
Set-TextValue "D12" "0.03069"
Set-TextValue "D13" "0.03054"
Set-TextValue "D14" "0.09323"
Set-TextValue "D15" "3.854"
Set-TextValue "D16" "0.001539"
Set-TextValue "D17" "0.04714"
Set-TextValue "D18" "0.0006049"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006176"
Set-TextValue "D20" "0.001266"
Set-TextValue "D21" "0.004591"
Set-TextValue "D22" "0.00008699"
Set-TextValue "D24" "2.179"
Set-TextValue "D40" "0.03777"
Set-TextValue "D41" "0.006309"
Set-TextValue "D43" "0.002605"
Set-TextValue "D44" "0.007970"
Set-TextValue "D45" "0.00005529"
Set-TextValue "D47" "0.5978"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue "D48" "0.01418"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"
